$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet1
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 155
$wsExhibit.Range("F3").Value = 1813
$wsExhibit.Range("F6").Value = 675
$wsExhibit.Range("F8").Value = 66
$wsExhibit.Range("F9").Value = 561
$wsExhibit.Range("F13").Value = 171
$wsExhibit.Range("F15").Value = 126
$wsExhibit.Range("F18").Value = 5181
$wsExhibit.Range("F19").Value = 61
$wsExhibit.Range("F21").Value = 123
$wsExhibit.Range("F22").Value = 2298
$wsExhibit.Range("F23").Value = 74
$wsExhibit.Range("F24").Value = 32
$wsExhibit.Range("F25").Value = 2141

# Sheet "全部类型" (All Types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 155
$wsAll.Range("F3").Value = 1813
$wsAll.Range("F6").Value = 675
$wsAll.Range("F8").Value = 66
$wsAll.Range("F9").Value = 561
$wsAll.Range("F13").Value = 171
$wsAll.Range("F15").Value = 126
$wsAll.Range("F18").Value = 5181
$wsAll.Range("F20").Value = 61
$wsAll.Range("F23").Value = 123
$wsAll.Range("F24").Value = 2298
$wsAll.Range("F25").Value = 74
$wsAll.Range("F27").Value = 32
$wsAll.Range("F28").Value = 2141
